$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 923-924, pushing the existing rows
# (923 onward) down by two -- matches the diff where row 923 becomes
# row 925, row 924 becomes row 926, ..., row 991 becomes row 993.
$ws.Rows("923:924").Insert()

# Fill the new row 923 with fresh data (Pintón quality).
$ws.Range("A923").Value = 5
$ws.Range("B923").Value = "Macroferia Regional de Talca"
$ws.Range("C923").Value = "Maule"
$ws.Range("D923").Value = 45021
$ws.Range("E923").Value = 7
$ws.Range("F923").Value = "Fruta"
$ws.Range("G923").Value = 100108
$ws.Range("H923").Value = "Tropicales y subtropicales"
$ws.Range("I923").Value = 100108006
$ws.Range("J923").Value = "Plátano"
$ws.Range("K923").Value = "Sin especificar"
$ws.Range("L923").Value = "Pintón"
$ws.Range("M923").Value = 850
$ws.Range("N923").Value = 20000
$ws.Range("O923").Value = 20000
$ws.Range("P923").Value = 20000
$ws.Range("Q923").Value = "`$/caja 20 kilos"
$ws.Range("R923").Value = "Ecuador"
$ws.Range("S923").Value = 1000
$ws.Range("T923").Value = 20

# Fill the new row 924 with fresh data (Primera Pintón quality).
$ws.Range("A924").Value = 5
$ws.Range("B924").Value = "Macroferia Regional de Talca"
$ws.Range("C924").Value = "Maule"
$ws.Range("D924").Value = 45021
$ws.Range("E924").Value = 7
$ws.Range("F924").Value = "Fruta"
$ws.Range("G924").Value = 100108
$ws.Range("H924").Value = "Tropicales y subtropicales"
$ws.Range("I924").Value = 100108006
$ws.Range("J924").Value = "Plátano"
$ws.Range("K924").Value = "Sin especificar"
$ws.Range("L924").Value = "Primera Pintón"
$ws.Range("M924").Value = 500
$ws.Range("N924").Value = 21000
$ws.Range("O924").Value = 21000
$ws.Range("P924").Value = 21000
$ws.Range("Q924").Value = "`$/caja 20 kilos"
$ws.Range("R924").Value = "Ecuador"
$ws.Range("S924").Value = 1050
$ws.Range("T924").Value = 20

# Give the new date cells the same date/time number format style as the
# rest of column D.
$ws.Range("D923:D924").NumberFormat = $ws.Range("D925").NumberFormat
